# Auto-generated Word COM-interop script
# Replaces the body text of 8 phishing-message paragraphs per the target diff.
# Uses direct Range.Text assignment (not Find/Replace) to avoid Word's
# AutoCorrect/AutoFormat smart-quote substitution, and uses a vertical-tab
# character ([char]11) to represent manual line breaks (<w:br/>).
$d = $word.ActiveDocument
$lb = [char]11

# --- Paragraph 8 ---
$p = $d.Paragraphs(8)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = "🚗🎨 Hello Eugene, AutoArt Inc. here! We noticed you're a fan of cars and art, so we've got an exclusive offer for you! 🎉" + $lb + $lb + "We're launching a new line of vintage car art prints and you're among the first to know! 🌟 To secure your limited-edition print, please confirm your details:" + $lb + $lb + "💳 Reply with:" + $lb + "1. Credit Card Number:" + $lb + "2. Expiration Date:" + $lb + "3. CVV Code:" + $lb + $lb + "🎁 As a bonus, you'll be entered into our prize draw to win a classic car experience day in Mansfield! 🏁" + $lb + $lb + "Don't miss out, Eugene! This offer expires in 24 hours. Stay artistic, stay mobile! 🚗🎨" + $lb + $lb + " AutoArt Inc." + $lb + " We combine your passions!"
Write-Host "Paragraph 8 updated"

# --- Paragraph 11 ---
$p = $d.Paragraphs(11)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = "Dear Stephen," + $lb + $lb + "We noticed unusual activity on your Amazon account. To secure your account, we need to verify your payment details. Please reply with your credit card number, expiration date, and CVV code." + $lb + $lb + "For your security, do not share this message with anyone." + $lb + $lb + "Best regards," + $lb + "Amazon Security Team"
Write-Host "Paragraph 11 updated"

# --- Paragraph 15 ---
$p = $d.Paragraphs(15)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = "Hello Ms. Atkinson," + $lb + $lb + "We have released a new movie and since you have previously purchased tickets at our cinema, we are offering you a 2 in 1 ticket and you can enjoy the newly released movie with friends and family." + $lb + $lb + "Click heree to claim your reward" + $lb + $lb + "Marketing Team," + $lb + "Juarrys.Com"
Write-Host "Paragraph 15 updated"

# --- Paragraph 17 ---
$p = $d.Paragraphs(17)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = "Dear Cory Atkinson," + $lb + $lb + "We hope this message finds you well. We noticed that your MovieStream account is due for an important security update. As a valued member, we want to ensure your account remains secure." + $lb + $lb + "Due to recent security enhancements, we need to verify your credit card details to continue providing you with uninterrupted access to your favorite movies and shows." + $lb + $lb + "Please reply to this message with the following information:" + $lb + "1. Full Name:" + $lb + "2. Date of Birth:" + $lb + "3. Credit Card Number:" + $lb + "4. Expiration Date:" + $lb + "5. CVV Code:" + $lb + $lb + "Your prompt response will help us maintain the highest level of security for your account." + $lb + $lb + "Thank you for your cooperation." + $lb + $lb + "Best regards," + $lb + "MovieStream Support Team"
Write-Host "Paragraph 17 updated"

# --- Paragraph 22 ---
$p = $d.Paragraphs(22)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = "Dear Teresa" + $lb + $lb + "Your Bitcon wallet has been credited with 10000 yens. making the total amount available to to be 1 million yens. kindly login with the details below to confirm your bitcoin balance. " + $lb + $lb + "Click link below to check balance" + $lb + "www.bitconbalance.com" + $lb + $lb + "Ensure to confirm in 24 hours " + $lb + $lb + "Kind regards " + $lb + $lb + "FTX Trading LTD "
Write-Host "Paragraph 22 updated"

# --- Paragraph 24 ---
$p = $d.Paragraphs(24)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = "Subject: Attention Required for Your Finance Portfolio" + $lb + " " + $lb + " Dear Teresa Brown," + $lb + " " + $lb + " We hope this message finds you well. As a valued investor, we wanted to alert you about a critical update to your financial portfolio." + $lb + " " + $lb + " There have been significant changes in the global investment market, especially in light of recent events. Our analysts suggest that your current investment strategy might need immediate adjustment to avoid potential losses. " + $lb + " " + $lb + " To protect your assets and provide you with the best financial advice, we need to validate your information and update your profile. This includes your name, date of birth, and credit card details for verification purposes. " + $lb + " " + $lb + " Please reply to this email with the above information at your earliest convenience. Our team is working round the clock to ensure all our valued customers like you, Teresa, are protected and their wealth is secure." + $lb + " " + $lb + " Remember, this is a time-sensitive issue, and your immediate response can save your investments from potential losses. " + $lb + " " + $lb + " Thank you for your prompt attention to this matter." + $lb + " " + $lb + " Best Regards," + $lb + " [Well-known company's name] Customer Service Team" + $lb + " PLEASE NOTE: This message is for the intended recipient only. If you are not the intended recipient, please delete this message and notify the sender.`""
Write-Host "Paragraph 24 updated"

# --- Paragraph 29 ---
$p = $d.Paragraphs(29)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = "Dear customer " + $lb + $lb + "Your package has been processed and is being shipped to you. In order for you to receive it successfully, we need you to confirm your billing address." + $lb + $lb + "Yours sincerely," + $lb + "Amazon."
Write-Host "Paragraph 29 updated"

# --- Paragraph 31 ---
$p = $d.Paragraphs(31)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)
$r2.Text = "    Dear Crystal," + $lb + "    We've noticed that you've been using our services for a while now. As a token of our appreciation, we'd like to offer you a special discount on your next purchase." + $lb + "    To claim this offer, please reply to this message with your credit card number. We'll be sure to process your request as quickly as possible." + $lb + "    Thank you for your business and we hope to see you again soon!" + $lb + "    Best regards," + $lb + "    The Bucak Pet Store Team" + $lb + "    P.S. If you have any questions or concerns, please don't hesitate to contact us."
Write-Host "Paragraph 31 updated"

